$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B22").Value = 20.5
$ws.Range("B23").Value = 717500
$ws.Range("B24").Value = 4026785.714285714
$ws.Range("B34").Value = 2894285.714285715
$ws.Range("B35").Value = 2894285.714285715
